# Add a "coming soon" placeholder message to the Charts tab and make it
# the active/selected sheet (matches the author's commit: "add coming soon
# message to xlsx templates").

$wb = $excel.ActiveWorkbook

$chartsSheet = $wb.Worksheets.Item("Charts")
$chartsSheet.Activate() | Out-Null
$chartsSheet.Range("A1").Value = "Automatically generated chart(s) coming soon to this tab."
$null = $chartsSheet.Range("A1").Select()
